# Generate Report for Handoff
#
# The localization job moved from "In Translation" to "Ready for handoff".
# Update the status + timestamp cells on all three sheets (Overview, zh-cn,
# de-de) and refresh the "latest handback" URL on the zh-cn sheet, then
# widen the Status columns so the new (longer) text still fits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns + HO xliff generate date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2017-02-09 14:00:46"

# ---------------------------------------------------------------------
# zh-cn sheet: Status + Latest Handoff Datetime + handback URL note
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2017-02-09 14:00:28"
$zhcn.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a072feab344118faf8bc8e6a6507da2c56498f26/e2e/c646d137-6169-4650-991f-2d337c5289f5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/617bcdaa826847eb3d81eb306c1505abdb0b6ce5/e2e/c646d137-6169-4650-991f-2d337c5289f5.md."

# ---------------------------------------------------------------------
# de-de sheet: Status
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Re-fit the Status columns now that "Ready for handoff" is longer than
# the previous "In Translation" text.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).AutoFit() | Out-Null
